$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Add graph data rows for poultry (266-269) and pigs (270-272).
# Column order of writes matters because it controls the order new strings
# are appended to the shared-string table, so we replicate the exact order
# used by the original author: all of column A, then column C, then column E
# for the poultry block, followed by the pigs block (with its slightly
# irregular A-column fill order).
# ---------------------------------------------------------------------------

# --- Column A (PyOrator variable names) -----------------------------------
$ws.Range("A266").Value = "poultry_n_excrete_nlim"
$ws.Range("A267").Value = "poultry_eggs_prod_nlim"
$ws.Range("A268").Value = "poultry_meat_prod_nlim"
$ws.Range("A269").Value = "poultry_manure_prod_nlim"

# --- Column C (short labels) ------------------------------------------------
$ws.Range("C266").Value = "N excreted by poultry (N Limited)"
$ws.Range("C267").Value = "Eggs produced by poultry (N Limited)"
$ws.Range("C268").Value = "Meat produced by poultry (N Limited)"
$ws.Range("C269").Value = "Manure produced by poultry (N Limited)"

# --- Column E (long descriptions) ------------------------------------------
$ws.Range("E266").Value = "Nitrogen excreted by poultry , crop growth limited by N availability"
$ws.Range("E267").Value = "Eggs produced by poultry, crop growth limited by N availability"
$ws.Range("E268").Value = "Meat produced by  poultry, crop growth limited by N availability"
$ws.Range("E269").Value = "Manure produced by poultry, crop growth limited by N availability"

# --- Pigs block --------------------------------------------------------------
$ws.Range("A270").Value = "pigs_n_excrete_nlim"

$ws.Range("C270").Value = "N excreted by pigs (N Limited)"
$ws.Range("C271").Value = "Meat produced by pigs (N Limited)"
$ws.Range("C272").Value = "Manure produced by pigs (N Limited)"

$ws.Range("E270").Value = "Nitrogen excreted by pigs, crop growth limited by N availability"
$ws.Range("E271").Value = "Meat produced by  pigs, crop growth limited by N availability"
$ws.Range("E272").Value = "Manure produced by pigs, crop growth limited by N availability"

$ws.Range("A272").Value = "pigs_manure_prod_nlim"
$ws.Range("A271").Value = "pigs_meat_prod_nlim"

# --- Column B (category -- reuses existing "livestock" string) -------------
$ws.Range("B266").Value = "livestock"
$ws.Range("B267").Value = "livestock"
$ws.Range("B268").Value = "livestock"
$ws.Range("B269").Value = "livestock"
$ws.Range("B270").Value = "livestock"
$ws.Range("B271").Value = "livestock"
$ws.Range("B272").Value = "livestock"

# --- Column F (units -- reuses existing "kg/y" string) ----------------------
$ws.Range("F266").Value = "kg/y"
$ws.Range("F267").Value = "kg/y"
$ws.Range("F268").Value = "kg/y"
$ws.Range("F269").Value = "kg/y"
$ws.Range("F270").Value = "kg/y"
$ws.Range("F271").Value = "kg/y"
$ws.Range("F272").Value = "kg/y"

# --- Column G (decimals flag -- reuses existing "2f" string) ----------------
$ws.Range("G266").Value = "2f"
$ws.Range("G267").Value = "2f"
$ws.Range("G268").Value = "2f"
$ws.Range("G269").Value = "2f"
$ws.Range("G270").Value = "2f"
$ws.Range("G271").Value = "2f"
$ws.Range("G272").Value = "2f"

# ---------------------------------------------------------------------------
# Row height - match the rest of the table (20.1 pt, custom height).
# ---------------------------------------------------------------------------
for ($r = 266; $r -le 272; $r++) {
    $ws.Rows.Item($r).RowHeight = 20.1
}

# ---------------------------------------------------------------------------
# Formatting to reproduce the existing cell styles used throughout the table.
# ---------------------------------------------------------------------------
foreach ($r in 266..272) {
    # Column B: horizontal center
    $ws.Range("B$r").HorizontalAlignment = -4108

    # Column F: horizontal + vertical center, wrap text
    $ws.Range("F$r").HorizontalAlignment = -4108
    $ws.Range("F$r").VerticalAlignment = -4108
    $ws.Range("F$r").WrapText = $true

    # Column G: horizontal + vertical center
    $ws.Range("G$r").HorizontalAlignment = -4108
    $ws.Range("G$r").VerticalAlignment = -4108
}

# Column E: rows 266, 267 and 270 use vertical-center + wrap text (style seen
# on the shorter of the long-text cells); 268, 269, 271, 272 keep plain
# default formatting, mirroring the source workbook exactly.
foreach ($r in 266, 267, 270) {
    $ws.Range("E$r").VerticalAlignment = -4108
    $ws.Range("E$r").WrapText = $true
}

# ---------------------------------------------------------------------------
# Update the active selection to match the end state of the edit (last cell
# touched by the author).
# ---------------------------------------------------------------------------
$ws.Range("A272").Select() | Out-Null
